$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value2 = "57.840.23"
$ws.Cells.Item(2,5).Value2 = "  -0.99%  "
$ws.Cells.Item(3,4).Value2 = "3.118.16"
$ws.Cells.Item(3,5).Value2 = "  -1.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Cells.Item(4,4).Value2 = "1.00"
$ws.Cells.Item(4,5).Value2 = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5,4).Value2 = "531.17"
$ws.Cells.Item(5,5).Value2 = "  -0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6,4).Value2 = "138.35"
$ws.Cells.Item(6,5).Value2 = "  -2.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Cells.Item(7,4).Value2 = "1.00"
$ws.Cells.Item(7,5).Value2 = "  -0.09%  "
$ws.Cells.Item(8,4).Value2 = "3.114.10"
$ws.Cells.Item(8,5).Value2 = "  -1.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Cells.Item(9,4).Value2 = "0.464"
$ws.Cells.Item(9,5).Value2 = "  +2.94%  "
$ws.Cells.Item(10,5).Value2 = "  +0.22%  "
$ws.Cells.Item(11,5).Value2 = "  -2.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Cells.Item(12,4).Value2 = "0.408"
$ws.Cells.Item(12,5).Value2 = "  +1.40%  "
$ws.Cells.Item(13,4).Value2 = "3.652.24"
$ws.Cells.Item(13,5).Value2 = "  -1.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Cells.Item(14,4).Value2 = "0.136"
$ws.Cells.Item(14,5).Value2 = "  +1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Cells.Item(15,4).Value2 = "25.47"
$ws.Cells.Item(15,5).Value2 = "  -1.54%  "
$ws.Cells.Item(16,5).Value2 = "  -2.47%  "
$ws.Cells.Item(17,4).Value2 = "57.844.57"
$ws.Cells.Item(17,5).Value2 = "  -1.21%  "
$ws.Cells.Item(18,4).Value2 = "3.121.06"
$ws.Cells.Item(18,5).Value2 = "  -2.04%  "
$ws.Cells.Item(19,5).Value2 = "  -3.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Cells.Item(20,4).Value2 = "12.62"
$ws.Cells.Item(20,5).Value2 = "  -2.71%  "
$ws.Cells.Item(21,5).Value2 = "  -1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Cells.Item(22,4).Value2 = "351.20"
$ws.Cells.Item(22,5).Value2 = "  -2.07%  "
$ws.Cells.Item(23,5).Value2 = "  +0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Cells.Item(24,4).Value2 = "68.92"
$ws.Cells.Item(24,5).Value2 = "  +0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Cells.Item(25,4).Value2 = "0.504"
$ws.Cells.Item(25,5).Value2 = "  -2.05%  "
$ws.Cells.Item(26,5).Value2 = "  -2.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Cells.Item(27,4).Value2 = "1.00"
$ws.Cells.Item(27,5).Value2 = "  +0.16%  "
$ws.Cells.Item(28,4).Value2 = "0.0₃0871"
$ws.Cells.Item(28,5).Value2 = "  -8.84%  "
$ws.Cells.Item(29,2).Value2 = "InternetComputer(DFINITY)"
$ws.Cells.Item(29,3).Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29,4).Value2 = "7.19"
$ws.Cells.Item(29,5).Value2 = "  -4.69%  "
$ws.Cells.Item(30,2).Value2 = "PancakeSwap"
$ws.Cells.Item(30,3).Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Cells.Item(30,4).Value2 = "1.87"
$ws.Cells.Item(30,5).Value2 = "  -2.35%  "
$ws.Cells.Item(31,2).Value2 = "RenderToken"
$ws.Cells.Item(31,3).Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").NumberFormat = "@"
$ws.Cells.Item(31,4).Value2 = "6.00"
$ws.Cells.Item(31,5).Value2 = "  -7.48%  "
$ws.Cells.Item(32,2).Value2 = "EthereumClassic"
$ws.Cells.Item(32,3).Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Cells.Item(32,4).Value2 = "21.26"
$ws.Cells.Item(32,5).Value2 = "  -0.79%  "
$ws.Cells.Item(33,2).Value2 = "NEARProtocol"
$ws.Cells.Item(33,3).Value2 = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Cells.Item(33,4).Value2 = "4.94"
$ws.Cells.Item(33,5).Value2 = "  +0.36%  "
$ws.Cells.Item(34,2).Value2 = "Fetch.AI"
$ws.Cells.Item(34,3).Value2 = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Cells.Item(34,4).Value2 = "1.14"
$ws.Cells.Item(34,5).Value2 = "  -6.71%  "
$ws.Cells.Item(35,2).Value2 = "Monero"
$ws.Cells.Item(35,3).Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Cells.Item(35,4).Value2 = "158.85"
$ws.Cells.Item(35,5).Value2 = "  +0.66%  "
$ws.Cells.Item(36,2).Value2 = "Aptos"
$ws.Cells.Item(36,3).Value2 = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Cells.Item(36,4).Value2 = "6.06"
$ws.Cells.Item(36,5).Value2 = "  -2.47%  "
$ws.Cells.Item(37,2).Value2 = "EnergySwap"
$ws.Cells.Item(37,3).Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").NumberFormat = "@"
$ws.Cells.Item(37,4).Value2 = "25.72"
$ws.Cells.Item(37,5).Value2 = "  -3.00%  "
$ws.Cells.Item(38,2).Value2 = "ImmutableX"
$ws.Cells.Item(38,3).Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Cells.Item(38,4).Value2 = "1.26"
$ws.Cells.Item(38,5).Value2 = "  -4.03%  "
$ws.Cells.Item(39,2).Value2 = "Stacks"
$ws.Cells.Item(39,3).Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Cells.Item(39,4).Value2 = "1.71"
$ws.Cells.Item(39,5).Value2 = "  +4.13%  "
$ws.Cells.Item(40,2).Value2 = "Hedera"
$ws.Cells.Item(40,3).Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Cells.Item(40,4).Value2 = "0.0668"
$ws.Cells.Item(40,5).Value2 = "  -1.24%  "
$ws.Cells.Item(41,2).Value2 = "Filecoin"
$ws.Cells.Item(41,3).Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Cells.Item(41,4).Value2 = "4.02"
$ws.Cells.Item(41,5).Value2 = "  -2.15%  "
$ws.Cells.Item(42,2).Value2 = "Mantle"
$ws.Cells.Item(42,3).Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42,4).Value2 = "0.697"
$ws.Cells.Item(42,5).Value2 = "  -1.89%  "
$ws.Cells.Item(43,2).Value2 = "Maker"
$ws.Cells.Item(43,3).Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(43,4).Value2 = "2.385.49"
$ws.Cells.Item(43,5).Value2 = "  +0.95%  "
$ws.Cells.Item(44,2).Value2 = "OKB"
$ws.Cells.Item(44,3).Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Cells.Item(44,4).Value2 = "36.97"
$ws.Cells.Item(44,5).Value2 = "  +0.19%  "
$ws.Cells.Item(45,2).Value2 = "FirstDigitalUSD"
$ws.Cells.Item(45,3).Value2 = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Cells.Item(45,4).Value2 = "1.00"
$ws.Cells.Item(45,5).Value2 = "  -0.25%  "
$ws.Cells.Item(46,2).Value2 = "RenzoRestakedETH"
$ws.Cells.Item(46,3).Value2 = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Cells.Item(46,4).Value2 = "3.154.80"
$ws.Cells.Item(46,5).Value2 = "  -1.95%  "
$ws.Cells.Item(47,2).Value2 = "VeChain"
$ws.Cells.Item(47,3).Value2 = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Cells.Item(47,4).Value2 = "0.0265"
$ws.Cells.Item(47,5).Value2 = "  -3.47%  "
$ws.Cells.Item(48,2).Value2 = "ONDO"
$ws.Cells.Item(48,3).Value2 = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Cells.Item(48,4).Value2 = "0.959"
$ws.Cells.Item(48,5).Value2 = "  -5.40%  "
$ws.Cells.Item(49,2).Value2 = "Cosmos"
$ws.Cells.Item(49,3).Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Cells.Item(49,4).Value2 = "6.04"
$ws.Cells.Item(49,5).Value2 = "  -0.85%  "
$ws.Cells.Item(50,2).Value2 = "InjectiveProtocol"
$ws.Cells.Item(50,3).Value2 = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Cells.Item(50,4).Value2 = "19.77"
$ws.Cells.Item(50,5).Value2 = "  -4.38%  "
$ws.Cells.Item(51,2).Value2 = "SuiNetwork"
$ws.Cells.Item(51,3).Value2 = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Cells.Item(51,4).Value2 = "0.742"
$ws.Cells.Item(51,5).Value2 = "  -3.30%  "
